$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing trial values (rows 2-4) with the new random results
$ws.Range("B2").Value = 9.6
$ws.Range("C2").Value = 5.4
$ws.Range("D2").Value = 6.4

$ws.Range("B3").Value = 0.8
$ws.Range("C3").Value = 7.8
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = 1.6
$ws.Range("C4").Value = 5.6
$ws.Range("D4").Value = 7.8

# Add the new "D (Groupe)" row of results
$ws.Range("A5").Value = "D (Groupe)"
$ws.Range("B5").Value = 2.2000000000000002
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 5

# Formulas for the new row mirror the ones above it (enter these while the
# row still carries default formatting, so the new cells don't inherit a
# stray number format from their precedents)
$ws.Range("E5").Formula = "=AVERAGE(B5:D5)"
$ws.Range("F5").Formula = "=STDEV(B5:D5)"
$ws.Range("G5").Formula = "=F5/E5"

# Now copy the row-above formatting down into the new row's cells
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null
$ws.Range("F4").Copy() | Out-Null
$ws.Range("F5").PasteSpecial(-4122) | Out-Null

# Match the final selection left behind in the workbook
$ws.Range("A5").Select() | Out-Null
